$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto the new
# header cells so they share the exact same style (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new columns I (I0) and J (IF), rows 2-21
$data = @(
    @(4, 6),
    @(5, 6),
    @(7, 8),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(3, 4),
    @(10, 10),
    @(8, 8),
    @(7, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(5, 5),
    @(4, 4),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
